$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 2.29025
$ws.Range("H2").Value = 6.870749999999999
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 4.166450999999999
$ws.Range("N2").Value = 12.499353
$ws.Range("O2").Value = 0.7126954333415383
$ws.Range("P2").Value = 0.7126954333415383
$ws.Range("Q2").Value = 9.542214402749998
$ws.Range("R2").Value = 85.87992962474999
$ws.Range("S2").Value = 0.7126954333415383
$ws.Range("T2").Value = 0.7126954333415383

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 2.29025
$ws.Range("H3").Value = 6.870749999999999
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.6655859999999999
$ws.Range("N3").Value = 1.996758
$ws.Range("O3").Value = 0.1138523176430159
$ws.Range("P3").Value = 0.1138523176430159
$ws.Range("Q3").Value = 1.5243583365
$ws.Range("R3").Value = 13.7192250285
$ws.Range("S3").Value = 0.1138523176430159
$ws.Range("T3").Value = 0.1138523176430159

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 2.29025
$ws.Range("H4").Value = 6.870749999999999
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 1.01401
$ws.Range("N4").Value = 3.04203
$ws.Range("O4").Value = 0.1734522490154458
$ws.Range("P4").Value = 0.1734522490154458
$ws.Range("Q4").Value = 2.3223364025
$ws.Range("R4").Value = 20.9010276225
$ws.Range("S4").Value = 0.1734522490154458
$ws.Range("T4").Value = 0.1734522490154458
